$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "IPs": rows 2-4 are re-synced against a fresher VirusTotal pull -
# the detected_url / detected_urls_positives / detected_urls_total /
# detected_urls_scan_date values shift down one row (row2->row3->row4) and
# row 2 receives the newest scan.
# ---------------------------------------------------------------------------
$ipsSheet = $wb.Worksheets.Item("IPs")

$ipsSheet.Range("G2").Value = "http://190.160.53.126/L3uZ0FJkzd00V/PnEJQO9BTZTIH75sat/"
$ipsSheet.Range("J2").Value = "2022-05-31 11:20:33"

$ipsSheet.Range("G3").Value = "http://190.160.53.126/cPzGZqbbrF2WtX/5aWaSd/XAWyg/AF8g1pcudGHa/5Tu4GPZZYIHhX7XZ5b7/"
$ipsSheet.Range("H3").Value = 10
$ipsSheet.Range("I3").Value = 94
$ipsSheet.Range("J3").Value = "2022-05-31 00:34:14"

$ipsSheet.Range("G4").Value = "http://190.160.53.126/"
$ipsSheet.Range("H4").Value = 8
$ipsSheet.Range("J4").Value = "2022-05-30 07:08:35"

# ---------------------------------------------------------------------------
# Sheet "URLs": a new "detected_urls_positives" column is inserted right
# after "detected_url" (old column D), pushing detected_urls_total,
# detected_urls_scan_date, Engine, Engine_detected and Engine_result one
# column to the right. The scan counters/date are refreshed for every row.
# ---------------------------------------------------------------------------
$urlsSheet = $wb.Worksheets.Item("URLs")

$urlsSheet.Range("D1").EntireColumn.Insert()

$urlsSheet.Range("D1").Value = "detected_urls_positives"

for ($row = 2; $row -le 7; $row++) {
    $urlsSheet.Cells.Item($row, 4).Value = 4
    $urlsSheet.Cells.Item($row, 5).Value = 86
    $urlsSheet.Cells.Item($row, 6).Value = "2022-05-31 10:11:18"
}
